# Update library jsoup error
# The "ainow" scraping config row is replaced with a "kyodonews" one, a
# couple of now-unused shared strings ("page/", "?s=") are dropped, and the
# selection on the URL sheet is moved to the (newly important) row 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("URL")

# --- Row 1 (headers) -------------------------------------------------
# Columns B:F were shifted over by the removal of the two header cells
# that used to sit in front of them; re-assert their labels explicitly.
$ws.Range("B1").Value = "CrawlSwitch"
$ws.Range("C1").Value = "SaveSwitch"
$ws.Range("D1").Value = "CrawlMethod"
$ws.Range("E1").Value = "SaveMethod"
$ws.Range("F1").Value = "PageNum"

$ws.Range("M1").Value = "click_Search"
$ws.Range("N1").Value = "class_Next_page"
$ws.Range("O1").Value = "Domain"

$ws.Range("U1").Value = "Item pages element"
$ws.Range("V1").Value = "備考(1ページ目の表示アイテム数)"

# --- Row 2 (kyodonews scraping config, replacing ainow) ---------------
$ws.Range("A2").Value = "kyodonews"
$ws.Range("B2").Value = 0
$ws.Range("G2").Value = "https://kyodonewsprwire.jp/search?s=(keyword)"
$ws.Range("H2").Value = "body > main > div > div > div > div > div > h4 > a"
$ws.Range("I2").Value = "body > main > div > div > div > div > div > h4 > a"
$ws.Range("J2").Value = "body > main > div > div > div > div > div > div > div > p.releae-up-date"
$ws.Range("K2").Value = "body > main > div > div > div > div.release-body > p:last-child"
$ws.Range("L2").Value = "body > main > div > div > div > div.release-body > p"
$ws.Range("N2").Value = "&page="
$ws.Range("O2").Value = "https://kyodonewsprwire.jp"
$ws.Range("P2").Value = "stp=1"
$ws.Range("Q2").Value = 25

# --- Selection moves to the whole of row 2 ----------------------------
$ws.Range("A2:XFD2").Select()
